# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" table (rows 16-33, cols B:G) is re-sorted from
# grouped-by-employee / period-descending to interleaved-by-employee /
# period-ascending (1911..2007), and the one-off 24292 Valor Mora
# (the partial-period figure) follows each employee oldest period (1911)
# instead of their most recent one (2007).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# B:Tipo Doc  C:N Doc  D:Nombre  E:Periodo Mora  F:Valor Mora  G:Salario Basico
$ws.Cells.Item(16, 2).Value = "CC"
$ws.Cells.Item(16, 3).Value = "1007208457"
$ws.Cells.Item(16, 4).Value = "MIGUEL ANTONIO ELLES BANGUERA"
$ws.Cells.Item(16, 5).Value = "1911"
$ws.Cells.Item(16, 6).Value = 33125
$ws.Cells.Item(16, 7).Value = 828116

$ws.Cells.Item(17, 2).Value = "CC"
$ws.Cells.Item(17, 3).Value = "1049927922"
$ws.Cells.Item(17, 4).Value = "WILMER DE JESUS MARTINEZ DE LOS REYES"
$ws.Cells.Item(17, 5).Value = "1911"
$ws.Cells.Item(17, 6).Value = 33125
$ws.Cells.Item(17, 7).Value = 828116

$ws.Cells.Item(18, 2).Value = "CC"
$ws.Cells.Item(18, 3).Value = "1007208457"
$ws.Cells.Item(18, 4).Value = "MIGUEL ANTONIO ELLES BANGUERA"
$ws.Cells.Item(18, 5).Value = "1912"
$ws.Cells.Item(18, 6).Value = 33125
$ws.Cells.Item(18, 7).Value = 828116

$ws.Cells.Item(19, 2).Value = "CC"
$ws.Cells.Item(19, 3).Value = "1049927922"
$ws.Cells.Item(19, 4).Value = "WILMER DE JESUS MARTINEZ DE LOS REYES"
$ws.Cells.Item(19, 5).Value = "1912"
$ws.Cells.Item(19, 6).Value = 33125
$ws.Cells.Item(19, 7).Value = 828116

$ws.Cells.Item(20, 2).Value = "CC"
$ws.Cells.Item(20, 3).Value = "1007208457"
$ws.Cells.Item(20, 4).Value = "MIGUEL ANTONIO ELLES BANGUERA"
$ws.Cells.Item(20, 5).Value = "2001"
$ws.Cells.Item(20, 6).Value = 33125
$ws.Cells.Item(20, 7).Value = 828116

$ws.Cells.Item(21, 2).Value = "CC"
$ws.Cells.Item(21, 3).Value = "1049927922"
$ws.Cells.Item(21, 4).Value = "WILMER DE JESUS MARTINEZ DE LOS REYES"
$ws.Cells.Item(21, 5).Value = "2001"
$ws.Cells.Item(21, 6).Value = 33125
$ws.Cells.Item(21, 7).Value = 828116

$ws.Cells.Item(22, 2).Value = "CC"
$ws.Cells.Item(22, 3).Value = "1007208457"
$ws.Cells.Item(22, 4).Value = "MIGUEL ANTONIO ELLES BANGUERA"
$ws.Cells.Item(22, 5).Value = "2002"
$ws.Cells.Item(22, 6).Value = 33125
$ws.Cells.Item(22, 7).Value = 828116

$ws.Cells.Item(23, 2).Value = "CC"
$ws.Cells.Item(23, 3).Value = "1049927922"
$ws.Cells.Item(23, 4).Value = "WILMER DE JESUS MARTINEZ DE LOS REYES"
$ws.Cells.Item(23, 5).Value = "2002"
$ws.Cells.Item(23, 6).Value = 33125
$ws.Cells.Item(23, 7).Value = 828116

$ws.Cells.Item(24, 2).Value = "CC"
$ws.Cells.Item(24, 3).Value = "1007208457"
$ws.Cells.Item(24, 4).Value = "MIGUEL ANTONIO ELLES BANGUERA"
$ws.Cells.Item(24, 5).Value = "2003"
$ws.Cells.Item(24, 6).Value = 33125
$ws.Cells.Item(24, 7).Value = 828116

$ws.Cells.Item(25, 2).Value = "CC"
$ws.Cells.Item(25, 3).Value = "1049927922"
$ws.Cells.Item(25, 4).Value = "WILMER DE JESUS MARTINEZ DE LOS REYES"
$ws.Cells.Item(25, 5).Value = "2003"
$ws.Cells.Item(25, 6).Value = 33125
$ws.Cells.Item(25, 7).Value = 828116

$ws.Cells.Item(26, 2).Value = "CC"
$ws.Cells.Item(26, 3).Value = "1007208457"
$ws.Cells.Item(26, 4).Value = "MIGUEL ANTONIO ELLES BANGUERA"
$ws.Cells.Item(26, 5).Value = "2004"
$ws.Cells.Item(26, 6).Value = 33125
$ws.Cells.Item(26, 7).Value = 828116

$ws.Cells.Item(27, 2).Value = "CC"
$ws.Cells.Item(27, 3).Value = "1049927922"
$ws.Cells.Item(27, 4).Value = "WILMER DE JESUS MARTINEZ DE LOS REYES"
$ws.Cells.Item(27, 5).Value = "2004"
$ws.Cells.Item(27, 6).Value = 33125
$ws.Cells.Item(27, 7).Value = 828116

$ws.Cells.Item(28, 2).Value = "CC"
$ws.Cells.Item(28, 3).Value = "1007208457"
$ws.Cells.Item(28, 4).Value = "MIGUEL ANTONIO ELLES BANGUERA"
$ws.Cells.Item(28, 5).Value = "2005"
$ws.Cells.Item(28, 6).Value = 33125
$ws.Cells.Item(28, 7).Value = 828116

$ws.Cells.Item(29, 2).Value = "CC"
$ws.Cells.Item(29, 3).Value = "1049927922"
$ws.Cells.Item(29, 4).Value = "WILMER DE JESUS MARTINEZ DE LOS REYES"
$ws.Cells.Item(29, 5).Value = "2005"
$ws.Cells.Item(29, 6).Value = 33125
$ws.Cells.Item(29, 7).Value = 828116

$ws.Cells.Item(30, 2).Value = "CC"
$ws.Cells.Item(30, 3).Value = "1007208457"
$ws.Cells.Item(30, 4).Value = "MIGUEL ANTONIO ELLES BANGUERA"
$ws.Cells.Item(30, 5).Value = "2006"
$ws.Cells.Item(30, 6).Value = 33125
$ws.Cells.Item(30, 7).Value = 828116

$ws.Cells.Item(31, 2).Value = "CC"
$ws.Cells.Item(31, 3).Value = "1049927922"
$ws.Cells.Item(31, 4).Value = "WILMER DE JESUS MARTINEZ DE LOS REYES"
$ws.Cells.Item(31, 5).Value = "2006"
$ws.Cells.Item(31, 6).Value = 33125
$ws.Cells.Item(31, 7).Value = 828116

$ws.Cells.Item(32, 2).Value = "CC"
$ws.Cells.Item(32, 3).Value = "1007208457"
$ws.Cells.Item(32, 4).Value = "MIGUEL ANTONIO ELLES BANGUERA"
$ws.Cells.Item(32, 5).Value = "2007"
$ws.Cells.Item(32, 6).Value = 24292
$ws.Cells.Item(32, 7).Value = 828116

$ws.Cells.Item(33, 2).Value = "CC"
$ws.Cells.Item(33, 3).Value = "1049927922"
$ws.Cells.Item(33, 4).Value = "WILMER DE JESUS MARTINEZ DE LOS REYES"
$ws.Cells.Item(33, 5).Value = "2007"
$ws.Cells.Item(33, 6).Value = 24292
$ws.Cells.Item(33, 7).Value = 828116
